$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readme")
$ws.Range("A1").Value = "test"
